$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FE)
$ws.Range("B2").Value = 1.45
$ws.Range("C2").Value = 0.37
$ws.Range("D2").Value = 1.31
$ws.Range("E2").Value = 1.01
$ws.Range("F2").Value = 0.4
$ws.Range("G2").Value = 1.52
$ws.Range("H2").Value = 0.38
$ws.Range("I2").Value = 0.07000000000000001
$ws.Range("J2").Value = 0.98
$ws.Range("K2").Value = 0.15

# Row 3 (FE+Disg)
$ws.Range("B3").Value = 0.98
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2.64
$ws.Range("E3").Value = 1.01
$ws.Range("F3").Value = 0.4
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 0.99
$ws.Range("K3").Value = 0.05

# Row 4 (FE+Disg+Var)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1.01
$ws.Range("F4").Value = 0.4
$ws.Range("G4").Value = 1.33
$ws.Range("H4").Value = 0.97
$ws.Range("I4").Value = 0.38
$ws.Range("J4").Value = 0.9
$ws.Range("K4").Value = 0.38
